# Add GitHub repository and Trello board URLs for the final project teams.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column E (bold, matching the other header cells)
$ws.Range("E1").Value = "Trello"
$ws.Range("E1").Font.Bold = $true

# Column D - GitHub repository links (D1 "Repository" already exists)
$ws.Range("D2").Value  = "https://github.com/480GameSuite/Arcade_App"
$ws.Range("D3").Value  = "https://github.com/leviwp48/Final-Project-Mobile-Apps"
$ws.Range("D4").Value  = "https://github.com/Rivey724/MobileAppsFinalProject"
$ws.Range("D5").Value  = "https://github.com/jaredconn/Greenthumbs"
$ws.Range("D6").Value  = "https://github.com/sam-alston/cs480-DTM"
$ws.Range("D7").Value  = "https://github.com/WillemTheWalrus/Dank_Spots"
$ws.Range("D8").Value  = "https://github.com/JackWitherell/ProjectAndroidSampler"
$ws.Range("D9").Value  = "https://github.com/zeefree/cyberfoxgame1"
$ws.Range("D10").Value = "https://github.com/KalimotxoGood/SmartVault"

# Column E - Trello board links (only for teams that have one)
$ws.Range("E3").Value  = "https://trello.com/b/7N4L97Pw/tactical-rpg"
$ws.Range("E4").Value  = "https://trello.com/b/jVBOEz9f/team-snes-final-project"
$ws.Range("E7").Value  = "https://trello.com/b/EWzPH2WB/database-set-up-and-connection"
$ws.Range("E10").Value = "https://trello.com/b/CLrzOFr5/smart-vault"

# Formatting touch-ups matching the authored view state
$ws.Columns.Item(1).ColumnWidth = 23.33203125
$ws.Columns.Item(2).ColumnWidth = 23.109375
$ws.Columns.Item(3).ColumnWidth = 20.33203125

$ws.Application.ActiveWindow.Zoom = 110
$ws.Range("E5").Select() | Out-Null

Write-Output "done"
